$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 4 data: D4 = "dim" (shared string already used in D3), E4 = 50
$ws.Range("D4").Value = "dim"
$ws.Range("E4").Value = 50

# Update the active selection from C8 to C7
$ws.Range("C7").Select()
